$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing "Вес" header from D2 to E2, and put the new "ГОСТ" header into D2
$ws.Range("E2").Value = "Вес"
$ws.Range("D2").Value = "ГОСТ"

# Copy the style of D2(old)/C2 to E2 so the new cell matches formatting (s="1")
$ws.Range("C2").Copy()
$ws.Range("E2").PasteSpecial(-4122) # xlPasteFormats

# Update the selected cell/range to I12
$ws.Range("I12").Select()
